$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7799339294433594
$ws.Range("B1").Value = 5.058628082275391
$ws.Range("C1").Value = 4.512070655822754
$ws.Range("D1").Value = 1.088398814201355
$ws.Range("E1").Value = 0.6696374416351318
